# Update need_to_buy.xlsx values for rows 2-7 (columns B, C, E, F)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = 8728.23697420251
    "C2" = 8371.41466040737
    "E2" = 4262.70212191531
    "F2" = -13.7451340698882

    "B3" = 8566.26435878059
    "C3" = 8348.38449624487
    "E3" = 3983.34252773214
    "F3" = 138.655292665709

    "B4" = 8465.18572246231
    "C4" = 7613.23514757351
    "E4" = 4563.79009970633
    "F4" = 132.209385303326

    "B5" = 2956.15781398948
    "C5" = 5157.44092140661
    "E5" = 4142.40915170159
    "F5" = 12.3270863795081

    "B6" = 2736.13408369351
    "C6" = 5017.73404703732
    "E6" = 4037.525883124
    "F6" = 2.13583042338837

    "B7" = 8780.1702078306
    "C7" = 8291.67789371877
    "E7" = 4910.42514252139
    "F7" = 174.92095984334
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
